# Remove the (redundant) explicit <w:contextualSpacing w:val="0"/> override
# from every paragraph's paragraph-properties (w:pPr). The element simply
# re-asserts the default ("off") behaviour, so dropping it is a pure,
# formatting-neutral cleanup.
#
# The Word object model exposed by this host does not surface a
# ParagraphFormat.ContextualSpacing property, so we go one level lower:
# read each paragraph's own WordOpenXML (a self-contained mini package
# wrapping just that paragraph), strip the <w:contextualSpacing/> element
# from the extracted <w:p>...</w:p> fragment, and feed the cleaned
# fragment back in via Range.InsertXML, which replaces the range's
# contents in place while leaving every sibling paragraph untouched.

$d = $word.ActiveDocument
$count = $d.Paragraphs.Count

for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    $rng = $para.Range
    $openXml = $rng.WordOpenXML

    if ($openXml -notmatch "contextualSpacing") {
        continue
    }

    $bodyStart = $openXml.IndexOf("<w:body>") + 8
    $bodyEnd = $openXml.IndexOf("</w:body>")
    $body = $openXml.Substring($bodyStart, $bodyEnd - $bodyStart)

    $pStart = $body.IndexOf("<w:p ")
    if ($pStart -lt 0) {
        $pStart = $body.IndexOf("<w:p>")
    }
    $pEnd = $body.IndexOf("</w:p>") + 6
    $fragment = $body.Substring($pStart, $pEnd - $pStart)

    $cleaned = $fragment -replace '<w:contextualSpacing[^>]*/>', ''

    $rng.InsertXML($cleaned)
}
